$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 125
$ws.Range("F3").Value = 323
$ws.Range("F4").Value = 188
$ws.Range("F5").Value = 1196
$ws.Range("F6").Value = 420
$ws.Range("F8").Value = 139
$ws.Range("F11").Value = 258
$ws.Range("F12").Value = 154
$ws.Range("F13").Value = 162
$ws.Range("F14").Value = 1395
$ws.Range("F15").Value = 528
$ws.Range("F16").Value = 203
$ws.Range("F17").Value = 316
$ws.Range("F19").Value = 735
$ws.Range("F20").Value = 1111
$ws.Range("F22").Value = 1923
$ws.Range("F23").Value = 2574
$ws.Range("F24").Value = 1331
$ws.Range("F25").Value = 58
$ws.Range("F26").Value = 256
$ws.Range("F27").Value = 381
$ws.Range("F28").Value = 998
$ws.Range("F29").Value = 783
$ws.Range("F30").Value = 1133
$ws.Range("F33").Value = 763
$ws.Range("F34").Value = 473
$ws.Range("F35").Value = 606
$ws.Range("F36").Value = 776
$ws.Range("F37").Value = 333
$ws.Range("F38").Value = 221

# sheet2 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 556

# sheet3 (Worksheets.Item(3))
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 863

# sheet4 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 863
$ws.Range("F6").Value = 125
$ws.Range("F7").Value = 323
$ws.Range("F8").Value = 188
$ws.Range("F11").Value = 1196
$ws.Range("F12").Value = 420
$ws.Range("F14").Value = 139
$ws.Range("F17").Value = 258
$ws.Range("F19").Value = 154
$ws.Range("F20").Value = 162
$ws.Range("F21").Value = 1395
$ws.Range("F22").Value = 528
$ws.Range("F23").Value = 203
$ws.Range("F24").Value = 316
$ws.Range("F26").Value = 1111
$ws.Range("F27").Value = 2574
$ws.Range("F29").Value = 1331
$ws.Range("F30").Value = 58
$ws.Range("F34").Value = 256
$ws.Range("F35").Value = 381
$ws.Range("F36").Value = 998
$ws.Range("F39").Value = 783
$ws.Range("F40").Value = 1133
$ws.Range("F41").Value = 763
$ws.Range("F42").Value = 473
$ws.Range("F43").Value = 606
$ws.Range("F44").Value = 776
$ws.Range("F45").Value = 333
$ws.Range("F48").Value = 221
